$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-111 is updated from serial date 45186 (2023-09-17)
# to serial date 45188 (2023-09-19).
for ($row = 2; $row -le 111; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
